$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (existing) - update Target cluster to "ECs" and refresh all numeric columns
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8330250000000001
$ws.Range("H2").Value = 2.499075
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.129292
$ws.Range("N2").Value = 0.387876
$ws.Range("O2").Value = 0.06615700391713267
$ws.Range("P2").Value = 0.06615700391713268
$ws.Range("Q2").Value = 0.1077034683
$ws.Range("R2").Value = 0.9693312147000002
$ws.Range("S2").Value = 0.06615700391713267
$ws.Range("T2").Value = 0.06615700391713268

# Row 3 (existing) - update Target cluster to "FAPs" and refresh all numeric columns
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8330250000000001
$ws.Range("H3").Value = 2.499075
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.4307096666666667
$ws.Range("N3").Value = 1.292129
$ws.Range("O3").Value = 0.220388431649395
$ws.Range("P3").Value = 0.220388431649395
$ws.Range("Q3").Value = 0.3587919200750001
$ws.Range("R3").Value = 3.229127280675001
$ws.Range("S3").Value = 0.220388431649395
$ws.Range("T3").Value = 0.220388431649395

# Row 4 (new row)
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt2"
$ws.Range("C4").Value = "Fzd3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8330250000000001
$ws.Range("H4").Value = 2.499075
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.394319
$ws.Range("N4").Value = 4.182957
$ws.Range("O4").Value = 0.7134545644334723
$ws.Range("P4").Value = 0.7134545644334724
$ws.Range("Q4").Value = 1.161502584975
$ws.Range("R4").Value = 10.453523264775
$ws.Range("S4").Value = 0.7134545644334723
$ws.Range("T4").Value = 0.7134545644334724
